# Daily attendance processing - 2026-02-07 04:14:59 UTC
# Reorders the "Recorded By" (column G) comma-separated name lists for
# several session rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = "Administrator, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Gehan Adel, Dr. Veronia Rafat"
$ws.Range("G3").Value  = "Administrator, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Veronia Rafat"
$ws.Range("G4").Value  = "Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Majorelle Magdy"
$ws.Range("G5").Value  = "Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Asmaa Reda"
$ws.Range("G6").Value  = "Dr. Menna tuâ€™Allah Medhat, Dr. Mohammad El-Tanany, Dr. Alshimaa Atef, Dr. Majorelle Magdy, Dr. Manar Montaser"
$ws.Range("G7").Value  = "Dr. Lamiaa Ossama, Dr. Menna tu'Alllah Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Nada Mohammad"
$ws.Range("G11").Value = "Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Range("G12").Value = "Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Amira Ibrahim, Dr. Dina Adel"
$ws.Range("G13").Value = "Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Amira Ibrahim"
$ws.Range("G17").Value = "Dr. Esraa Samy, Dr. Mohammad Safwat"
$ws.Range("G19").Value = "Dr. Rania Ahmad Youssef, Dr. Mariam Toma Gerges"
$ws.Range("G30").Value = "Dr. Shorok Mohammad, Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Wafaa Ebida"
